{"js": "// Update the date heading and the 25 multiplication problems in the table\n// (5 populated rows x 5 columns out of 20 rows total; the other rows are\n// blank spacer rows). Each cell is addressed by its row/col index and its\n// text is replaced in place via a Range so run/paragraph formatting is\n// preserved (only the <w:t> content changes, matching the source diff).\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Heading paragraph: \"2024-09-23 Monday\" -> \"2024-09-24 Tuesday\"\nparagraphs.items[0].getRange().insertText(\"2024-09-24 Tuesday\", Word.InsertLocation.replace);\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Row index (0-based) -> new values for the 5 cells in that row.\nconst rowUpdates = {\n  0: [\"58\u00d742=\", \"46\u00d789=\", \"59\u00d799=\", \"50\u00d752=\", \"79\u00d769=\"],\n  4: [\"68\u00d735=\", \"17\u00d767=\", \"60\u00d714=\", \"76\u00d756=\", \"35\u00d745=\"],\n  9: [\"88\u00d730=\", \"21\u00d757=\", \"77\u00d730=\", \"99\u00d749=\", \"63\u00d765=\"],\n  14: [\"29\u00d713=\", \"35\u00d774=\", \"43\u00d723=\", \"50\u00d735=\", \"33\u00d799=\"],\n  19: [\"53\u00d764=\", \"30\u00d762=\", \"44\u00d726=\", \"82\u00d780=\", \"43\u00d737=\"],\n};\n\nfor (const rowIndex of Object.keys(rowUpdates)) {\n  const newValues = rowUpdates[rowIndex];\n  for (let col = 0; col < newValues.length; col++) {\n    const cell = table.getCell(parseInt(rowIndex, 10), col);\n    cell.getRange().insertText(newValues[col], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 multiplication problems in the table.\n# The table has 20 rows x 5 columns; only rows 1, 5, 10, 15, 20 (1-based)\n# contain text, the rest are blank spacer rows. Each cell's Range.Text is\n# replaced in place so existing run/paragraph formatting is preserved.\n\n$d = $word.ActiveDocument\n\n$p = $d.Paragraphs.Item(1)\n$p.Range.Text = \"2024-09-24 Tuesday\"\n\n$t = $d.Tables.Item(1)\n\n$rowUpdates = @{\n    1  = @(\"58\u00d742=\", \"46\u00d789=\", \"59\u00d799=\", \"50\u00d752=\", \"79\u00d769=\")\n    5  = @(\"68\u00d735=\", \"17\u00d767=\", \"60\u00d714=\", \"76\u00d756=\", \"35\u00d745=\")\n    10 = @(\"88\u00d730=\", \"21\u00d757=\", \"77\u00d730=\", \"99\u00d749=\", \"63\u00d765=\")\n    15 = @(\"29\u00d713=\", \"35\u00d774=\", \"43\u00d723=\", \"50\u00d735=\", \"33\u00d799=\")\n    20 = @(\"53\u00d764=\", \"30\u00d762=\", \"44\u00d726=\", \"82\u00d780=\", \"43\u00d737=\")\n}\n\nforeach ($rowIndex in $rowUpdates.Keys) {\n    $newValues = $rowUpdates[$rowIndex]\n    for ($col = 1; $col -le $newValues.Length; $col++) {\n        $cell = $t.Cell($rowIndex, $col)\n        $cell.Range.Text = $newValues[$col - 1]\n    }\n}\n\nWrite-Output \"done\"\n"}
